$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.228
$ws.Range("C6").Value = -11.864
$ws.Range("C7").Value = -12.793
$ws.Range("D7").Value = -7.518000000000001
$ws.Range("D12").Value = -7.562
$ws.Range("D15").Value = -8.423
$ws.Range("C16").Value = -13.312
$ws.Range("C20").Value = -12.047
$ws.Range("D20").Value = -7.994
$ws.Range("D21").Value = -7.967000000000001
$ws.Range("D22").Value = -7.906000000000001
$ws.Range("D23").Value = -7.997
$ws.Range("C28").Value = -12.848
$ws.Range("C29").Value = -12.035
$ws.Range("D29").Value = -7.325
$ws.Range("C32").Value = -13.222
$ws.Range("D34").Value = -7.936999999999999
$ws.Range("C40").Value = -12.151
$ws.Range("D42").Value = -8.107000000000001
$ws.Range("D43").Value = -7.795
$ws.Range("D44").Value = -7.580999999999999
$ws.Range("D45").Value = -7.498
$ws.Range("C46").Value = -13.549
$ws.Range("D46").Value = -8.605
$ws.Range("D50").Value = -8.178000000000001
$ws.Range("C51").Value = -12.199
$ws.Range("D51").Value = -7.598999999999999
$ws.Range("C52").Value = -11.614
$ws.Range("C57").Value = -13.844
$ws.Range("C59").Value = -12.718
$ws.Range("C62").Value = -13.737
$ws.Range("C66").Value = -11.038
$ws.Range("D66").Value = -7.683000000000002
$ws.Range("D67").Value = -7.194
$ws.Range("C73").Value = -12.401
$ws.Range("C74").Value = -12.144
$ws.Range("D79").Value = -7.465000000000001
$ws.Range("D84").Value = -8.266
$ws.Range("C92").Value = -11.037
$ws.Range("D92").Value = -6.772
$ws.Range("D97").Value = -8.312999999999999
$ws.Range("C100").Value = -12.706
